$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert a new row at 36 (pushing the existing rows 36-44
# down to 37-45) and populate it with this week's "Ajo" price data.
$ws.Rows("36:36").Insert()

$ws.Range("A36").Value = 1
$ws.Range("B36").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C36").Value = "Arica y Parinacota"
$ws.Range("D36").Value = 45093
$ws.Range("E36").Value = 15
$ws.Range("F36").Value = 100112003
$ws.Range("G36").Value = "Ajo"
$ws.Range("H36").Value = "Chino"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 400
$ws.Range("K36").Value = 16000
$ws.Range("L36").Value = 18000
$ws.Range("M36").Value = 17000
$ws.Range("N36").Value = '$/caja 10 kilos'
$ws.Range("O36").Value = "China"
$ws.Range("P36").Value = 1700
$ws.Range("Q36").Value = 10
$ws.Range("R36").Value = "Hortaliza"
